$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.562128560411576
$ws.Range("C2").Value = 1.809441918582433
$ws.Range("D2").Value = 1.552136229771481
$ws.Range("E2").Value = 1.414509557120016
$ws.Range("B3").Value = 2.576937016845974
$ws.Range("C3").Value = 1.820836472639407
$ws.Range("D3").Value = 1.563193023149681
$ws.Range("E3").Value = 1.425213289967862
$ws.Range("B4").Value = 2.526949918108024
$ws.Range("C4").Value = 1.782668010042132
$ws.Range("D4").Value = 1.526684029639687
$ws.Range("E4").Value = 1.390138974964715
$ws.Range("B5").Value = 2.571187619835253
$ws.Range("C5").Value = 1.817028306682256
$ws.Range("D5").Value = 1.545217619614819
$ws.Range("E5").Value = 1.423685679465119
$ws.Range("B6").Value = 2.580933269305367
$ws.Range("C6").Value = 1.823976775771101
$ws.Range("D6").Value = 1.551005964670673
$ws.Range("E6").Value = 1.428914950493217
$ws.Range("B7").Value = 2.555471140536096
$ws.Range("C7").Value = 1.804522491512453
$ws.Range("D7").Value = 1.549486314736102
$ws.Range("E7").Value = 1.413520031540223
$ws.Range("B8").Value = 2.568469745368881
$ws.Range("C8").Value = 1.815118152670409
$ws.Range("D8").Value = 1.573590969358645
$ws.Range("E8").Value = 1.419552146857844
$ws.Range("B9").Value = 2.57466267579235
$ws.Range("C9").Value = 1.820110013471357
$ws.Range("D9").Value = 1.563634562990191
$ws.Range("E9").Value = 1.426263507133462
$ws.Range("B10").Value = 2.253757000171591
$ws.Range("C10").Value = 1.57578486270646
$ws.Range("D10").Value = 1.344117272470109
$ws.Range("E10").Value = 1.225441711943755
$ws.Range("B11").Value = 2.186605345224173
$ws.Range("C11").Value = 1.524747233310899
$ws.Range("D11").Value = 1.298091163975458
$ws.Range("E11").Value = 1.183188477810487
$ws.Range("B12").Value = 1.858851259405834
$ws.Range("C12").Value = 1.275774180011838
$ws.Range("D12").Value = 1.073177682413108
$ws.Range("E12").Value = 0.9765449706184237
$ws.Range("B13").Value = 2.23224328459987
$ws.Range("C13").Value = 1.559307860733568
$ws.Range("D13").Value = 1.329072515563819
$ws.Range("E13").Value = 1.21168381432692
